$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "land" column right after "attraction" (current column B),
# shifting the existing "lat"/"lon" columns from B/C to C/D.
$ws.Columns("B").Insert()
$ws.Columns("B").ColumnWidth = 36.5

# Header
$ws.Range("B1").Value = "land"

# Land values for each attraction row (2-12)
$ws.Range("B2").Value  = "Illumination's Minion Land"
$ws.Range("B3").Value  = "Production Central"
$ws.Range("B4").Value  = "The Wizarding World of Harry Potter"
$ws.Range("B5").Value  = "New York"
$ws.Range("B6").Value  = "The Wizarding World of Harry Potter"
$ws.Range("B7").Value  = "World Expo"
$ws.Range("B8").Value  = "World Expo"
$ws.Range("B9").Value  = "New York"
$ws.Range("B10").Value = "Production Central"
$ws.Range("B11").Value = "Wood Woodpecker's Kidzone"
$ws.Range("B12").Value = "San Francisco"

# The geo hyperlink originally anchored on the lon cell (old C2) needs to
# move to the new lon column (D2) now that a column was inserted.
# Adding the hyperlink with a display string overwrites the cell's value and
# applies the built-in "Hyperlink" style, so restore the original numeric
# longitude value and clear the style afterwards to match the source data.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "geo:28.475272,-81.468103", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "geo:28.475272,-81.468103")
$ws.Range("D2").Value = -81.468102999999999
$ws.Range("D2").Style = "Normal"
# Adding a hyperlink registers a built-in "Hyperlink" cell style; the source
# workbook never used it (the link was added at the XML level), so drop it.
$wb.Styles.Item("Hyperlink").Delete()

# Reflect the selection state recorded in the saved workbook.
[void]$ws.Range("E1:E1048576").Select()
